# Add "Request/sec" column (L = H*60, N = K*60) to each of the 4 report
# sheets (V1..V4), matching the jmeter report's existing "Throughput"
# columns (H on the summary table, K on the aggregate table).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Header cells for the new "Request/sec" column.
    $ws.Range("L2").Value = "Request/sec"
    $ws.Range("N9").Value = "Request/sec"

    # Data rows: Request/sec = Throughput (requests/sec) * 60 (sec/min).
    $ws.Range("L3").Formula = "=H3*60"
    $ws.Range("L4").Formula = "=H4*60"
    $ws.Range("N10").Formula = "=K10*60"
    $ws.Range("N11").Formula = "=K11*60"

    # Move the sheet's active selection onto the newly added column.
    $ws.Range("L2:L4").Select()
}

$wb.Worksheets.Item("V1").Activate()
